$p = $ppt.ActivePresentation
$ds = $p.Designs
$d1 = $ds.Item(1)
$d1.Name = "TestName"
Write-Output "set done"
